$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sire/dam columns (F and G) with a single "sex" value column (D) = "U"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"      # D: sex = "U"
    $ws.Cells.Item($r, 6).Value = $null    # F: sire -> clear
    $ws.Cells.Item($r, 7).Value = $null    # G: dam -> clear
}

$ws.Range("D7").Select()
